$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 18.91491841113828
$ws.Cells.Item(2, 3).Value = 8.934646501143071
$ws.Cells.Item(2, 4).Value = 8.255381932655309
$ws.Cells.Item(2, 5).Value = 12.38757432589048
$ws.Cells.Item(2, 6).Value = 33.69482902560704
$ws.Cells.Item(2, 8).Value = 7.344005520526261
$ws.Cells.Item(2, 9).Value = 24.77953575933526
$ws.Cells.Item(2, 10).Value = 9.623628999666686
$ws.Cells.Item(2, 12).Value = 11.53494986194546
$ws.Cells.Item(2, 15).Value = 25.83468950717293
$ws.Cells.Item(3, 2).Value = 18.31392817680905
$ws.Cells.Item(3, 3).Value = 8.553639185446569
$ws.Cells.Item(3, 4).Value = 8.248553358349813
$ws.Cells.Item(3, 5).Value = 12.42236024627984
$ws.Cells.Item(3, 6).Value = 33.82991092919632
$ws.Cells.Item(3, 8).Value = 7.344005520526261
$ws.Cells.Item(3, 9).Value = 24.93851085218981
$ws.Cells.Item(3, 10).Value = 9.648996252664226
$ws.Cells.Item(3, 12).Value = 11.50412389053931
$ws.Cells.Item(3, 15).Value = 25.95752624055027
$ws.Cells.Item(4, 2).Value = 17.93549525465857
$ws.Cells.Item(4, 3).Value = 8.30971540045657
$ws.Cells.Item(4, 4).Value = 8.245132061207565
$ws.Cells.Item(4, 5).Value = 12.4451525064107
$ws.Cells.Item(4, 6).Value = 33.9218742675852
$ws.Cells.Item(4, 8).Value = 7.344005520526261
$ws.Cells.Item(4, 9).Value = 25.04180605976789
$ws.Cells.Item(4, 10).Value = 9.665372723096111
$ws.Cells.Item(4, 12).Value = 11.48640380535896
$ws.Cells.Item(4, 15).Value = 26.03951817855526
$ws.Cells.Item(5, 2).Value = 17.77912295581072
$ws.Cells.Item(5, 3).Value = 8.207902639933264
$ws.Cells.Item(5, 4).Value = 8.243933060137239
$ws.Cells.Item(5, 5).Value = 12.45480160939264
$ws.Cells.Item(5, 6).Value = 33.96161145667689
$ws.Cells.Item(5, 8).Value = 7.344005520526261
$ws.Cells.Item(5, 9).Value = 25.08532968212033
$ws.Cells.Item(5, 10).Value = 9.672248260332472
$ws.Cells.Item(5, 12).Value = 11.47949057106102
$ws.Cells.Item(5, 15).Value = 26.07457817576494
$ws.Cells.Item(6, 2).Value = 17.75303366025093
$ws.Cells.Item(6, 3).Value = 8.190854204261591
$ws.Cells.Item(6, 4).Value = 8.243745791017277
$ws.Cells.Item(6, 5).Value = 12.45642566087164
$ws.Cells.Item(6, 6).Value = 33.96834614808441
$ws.Cells.Item(6, 8).Value = 7.344005520526261
$ws.Cells.Item(6, 9).Value = 25.09264313726559
$ws.Cells.Item(6, 10).Value = 9.673402155735605
$ws.Cells.Item(6, 12).Value = 11.47836134776432
$ws.Cells.Item(6, 15).Value = 26.08049924515172
$ws.Cells.Item(7, 2).Value = 17.93339480742511
$ws.Cells.Item(7, 3).Value = 8.308351942563956
$ws.Cells.Item(7, 4).Value = 8.245115099118925
$ws.Cells.Item(7, 5).Value = 12.44528117483152
$ws.Cells.Item(7, 6).Value = 33.92240103093035
$ws.Cells.Item(7, 8).Value = 7.344005520526261
$ws.Cells.Item(7, 9).Value = 25.04238724358961
$ws.Cells.Item(7, 10).Value = 9.665464630331805
$ws.Cells.Item(7, 12).Value = 11.48630931911156
$ws.Cells.Item(7, 15).Value = 26.03998434463245
$ws.Cells.Item(8, 2).Value = 18.70977401060528
$ws.Cells.Item(8, 3).Value = 8.805402187622301
$ws.Cells.Item(8, 4).Value = 8.252868063862104
$ws.Cells.Item(8, 5).Value = 12.39927133039408
$ws.Cells.Item(8, 6).Value = 33.73952826208895
$ws.Cells.Item(8, 8).Value = 7.344005520526261
$ws.Cells.Item(8, 9).Value = 24.83317111780618
$ws.Cells.Item(8, 10).Value = 9.632209772939254
$ws.Cells.Item(8, 12).Value = 11.52407314941284
$ws.Cells.Item(8, 15).Value = 25.8756774466936
$ws.Cells.Item(9, 2).Value = 20.14910027797533
$ws.Cells.Item(9, 3).Value = 9.697246936767856
$ws.Cells.Item(9, 4).Value = 8.274138136974305
$ws.Cells.Item(9, 5).Value = 12.32039421649977
$ws.Cells.Item(9, 6).Value = 33.452826637689
$ws.Cells.Item(9, 8).Value = 7.344005520526261
$ws.Cells.Item(9, 9).Value = 24.46796438094947
$ws.Cells.Item(9, 10).Value = 9.573323843198493
$ws.Cells.Item(9, 12).Value = 11.60751024282115
$ws.Cells.Item(9, 15).Value = 25.60579458203549
$ws.Cells.Item(10, 2).Value = 21.14589348913665
$ws.Cells.Item(10, 3).Value = 10.29791472249729
$ws.Cells.Item(10, 4).Value = 8.293387824673225
$ws.Cells.Item(10, 5).Value = 12.26932443426507
$ws.Cells.Item(10, 6).Value = 33.28645215165197
$ws.Cells.Item(10, 8).Value = 7.344005520526261
$ws.Cells.Item(10, 9).Value = 24.22707387350372
$ws.Cells.Item(10, 10).Value = 9.533878090855417
$ws.Cells.Item(10, 12).Value = 11.67425736422484
$ws.Cells.Item(10, 15).Value = 25.4396738810175
$ws.Cells.Item(11, 2).Value = 21.58441245985476
$ws.Cells.Item(11, 3).Value = 10.55867527605919
$ws.Cells.Item(11, 4).Value = 8.302914214878539
$ws.Cells.Item(11, 5).Value = 12.24757791702846
$ws.Cells.Item(11, 6).Value = 33.22046178794712
$ws.Cells.Item(11, 8).Value = 7.344005520526261
$ws.Cells.Item(11, 9).Value = 24.12343412451679
$ws.Cells.Item(11, 10).Value = 9.516753825840265
$ws.Cells.Item(11, 12).Value = 11.70574246609774
$ws.Cells.Item(11, 15).Value = 25.37114145514664
$ws.Cells.Item(12, 2).Value = 21.74819509582455
$ws.Cells.Item(12, 3).Value = 10.65558005546342
$ws.Cells.Item(12, 4).Value = 8.306630672125367
$ws.Cells.Item(12, 5).Value = 12.23955609640388
$ws.Cells.Item(12, 6).Value = 33.19687360181638
$ws.Cells.Item(12, 8).Value = 7.344005520526261
$ws.Cells.Item(12, 9).Value = 24.08504278387185
$ws.Cells.Item(12, 10).Value = 9.510386571330141
$ws.Cells.Item(12, 12).Value = 11.71782070112406
$ws.Cells.Item(12, 15).Value = 25.34620657256538
$ws.Cells.Item(13, 2).Value = 21.71302464548928
$ws.Cells.Item(13, 3).Value = 10.63479231789674
$ws.Cells.Item(13, 4).Value = 8.305825444816689
$ws.Cells.Item(13, 5).Value = 12.24127426834506
$ws.Cells.Item(13, 6).Value = 33.20189132640731
$ws.Cells.Item(13, 8).Value = 7.344005520526261
$ws.Cells.Item(13, 9).Value = 24.09327302402257
$ws.Cells.Item(13, 10).Value = 9.51175266387312
$ws.Cells.Item(13, 12).Value = 11.71521261431787
$ws.Cells.Item(13, 15).Value = 25.35153144160323
$ws.Cells.Item(14, 2).Value = 21.59793323374789
$ws.Cells.Item(14, 3).Value = 10.5666847850994
$ws.Cells.Item(14, 4).Value = 8.303217794804791
$ws.Cells.Item(14, 5).Value = 12.24691368864819
$ws.Cells.Item(14, 6).Value = 33.21849305850456
$ws.Cells.Item(14, 8).Value = 7.344005520526261
$ws.Cells.Item(14, 9).Value = 24.1202585128345
$ws.Cells.Item(14, 10).Value = 9.516227639609649
$ws.Cells.Item(14, 12).Value = 11.70673306672036
$ws.Cells.Item(14, 15).Value = 25.3690696395526
$ws.Cells.Item(15, 2).Value = 21.527136588841
$ws.Cells.Item(15, 3).Value = 10.52472620372478
$ws.Cells.Item(15, 4).Value = 8.301634681454638
$ws.Cells.Item(15, 5).Value = 12.25039573586874
$ws.Cells.Item(15, 6).Value = 33.22884474098324
$ws.Cells.Item(15, 8).Value = 7.344005520526261
$ws.Cells.Item(15, 9).Value = 24.13689923136718
$ws.Cells.Item(15, 10).Value = 9.518983955311782
$ws.Cells.Item(15, 12).Value = 11.70155917705193
$ws.Cells.Item(15, 15).Value = 25.37994486464943
$ws.Cells.Item(16, 2).Value = 21.11692354978569
$ws.Cells.Item(16, 3).Value = 10.28061821397764
$ws.Cells.Item(16, 4).Value = 8.292780599729193
$ws.Cells.Item(16, 5).Value = 12.27077546859626
$ws.Cells.Item(16, 6).Value = 33.29096044843551
$ws.Cells.Item(16, 8).Value = 7.344005520526261
$ws.Cells.Item(16, 9).Value = 24.23396657297814
$ws.Cells.Item(16, 10).Value = 9.535013650971047
$ws.Cells.Item(16, 12).Value = 11.67222181789647
$ws.Cells.Item(16, 15).Value = 25.44429466035603
$ws.Cells.Item(17, 2).Value = 20.86135177960088
$ws.Cells.Item(17, 3).Value = 10.12763408288378
$ws.Cells.Item(17, 4).Value = 8.287544870600696
$ws.Cells.Item(17, 5).Value = 12.28365787117873
$ws.Cells.Item(17, 6).Value = 33.33155437197452
$ws.Cells.Item(17, 8).Value = 7.344005520526261
$ws.Cells.Item(17, 9).Value = 24.2950364149105
$ws.Cells.Item(17, 10).Value = 9.545056926354301
$ws.Cells.Item(17, 12).Value = 11.65450750014817
$ws.Cells.Item(17, 15).Value = 25.48557701027178
$ws.Cells.Item(18, 2).Value = 20.71295533514545
$ws.Cells.Item(18, 3).Value = 10.03846851759511
$ws.Cells.Item(18, 4).Value = 8.284605959948662
$ws.Cells.Item(18, 5).Value = 12.29120732814401
$ws.Cells.Item(18, 6).Value = 33.35581499631778
$ws.Cells.Item(18, 8).Value = 7.344005520526261
$ws.Cells.Item(18, 9).Value = 24.33072149140852
$ws.Cells.Item(18, 10).Value = 9.550910749541078
$ws.Cells.Item(18, 12).Value = 11.64442461047141
$ws.Cells.Item(18, 15).Value = 25.5099835590044
$ws.Cells.Item(19, 2).Value = 20.66247484307321
$ws.Cells.Item(19, 3).Value = 10.0080785674701
$ws.Cells.Item(19, 4).Value = 8.2836234036709
$ws.Cells.Item(19, 5).Value = 12.29378747697715
$ws.Cells.Item(19, 6).Value = 33.36418566020669
$ws.Cells.Item(19, 8).Value = 7.344005520526261
$ws.Cells.Item(19, 9).Value = 24.34289990296446
$ws.Cells.Item(19, 10).Value = 9.552906029073828
$ws.Cells.Item(19, 12).Value = 11.64102908127459
$ws.Cells.Item(19, 15).Value = 25.51836074916329
$ws.Cells.Item(20, 2).Value = 20.88870352914701
$ws.Cells.Item(20, 3).Value = 10.14404126307694
$ws.Cells.Item(20, 4).Value = 8.288094726798747
$ws.Cells.Item(20, 5).Value = 12.282272048064
$ws.Cells.Item(20, 6).Value = 33.32713863217814
$ws.Cells.Item(20, 8).Value = 7.344005520526261
$ws.Cells.Item(20, 9).Value = 24.28847753712081
$ws.Cells.Item(20, 10).Value = 9.54397981698761
$ws.Cells.Item(20, 12).Value = 11.65638229863064
$ws.Cells.Item(20, 15).Value = 25.4811138794187
$ws.Cells.Item(21, 2).Value = 21.631801064858
$ws.Cells.Item(21, 3).Value = 10.58673985928062
$ws.Cells.Item(21, 4).Value = 8.303980779516142
$ws.Cells.Item(21, 5).Value = 12.24525147418075
$ws.Cells.Item(21, 6).Value = 33.21357865302384
$ws.Cells.Item(21, 8).Value = 7.344005520526261
$ws.Cells.Item(21, 9).Value = 24.11230901947834
$ws.Cells.Item(21, 10).Value = 9.514910050770844
$ws.Cells.Item(21, 12).Value = 11.70921954166615
$ws.Cells.Item(21, 15).Value = 25.36389060886549
$ws.Cells.Item(22, 2).Value = 22.10414843704795
$ws.Cells.Item(22, 3).Value = 10.86533099663499
$ws.Cells.Item(22, 4).Value = 8.314997874977479
$ws.Cells.Item(22, 5).Value = 12.22229835810298
$ws.Cells.Item(22, 6).Value = 33.14752931260676
$ws.Cells.Item(22, 8).Value = 7.344005520526261
$ws.Cells.Item(22, 9).Value = 24.00215523971771
$ws.Cells.Item(22, 10).Value = 9.49659493059257
$ws.Cells.Item(22, 12).Value = 11.74465537468627
$ws.Cells.Item(22, 15).Value = 25.29320774303746
$ws.Cells.Item(23, 2).Value = 21.85330431492731
$ws.Cells.Item(23, 3).Value = 10.71763690784692
$ws.Cells.Item(23, 4).Value = 8.309060335511633
$ws.Cells.Item(23, 5).Value = 12.23443538571989
$ws.Cells.Item(23, 6).Value = 33.18203145883851
$ws.Cells.Item(23, 8).Value = 7.344005520526261
$ws.Cells.Item(23, 9).Value = 24.06049036950441
$ws.Cells.Item(23, 10).Value = 9.506307682613393
$ws.Cells.Item(23, 12).Value = 11.72566186806246
$ws.Cells.Item(23, 15).Value = 25.33038834964425
$ws.Cells.Item(24, 2).Value = 20.87634235589023
$ws.Cells.Item(24, 3).Value = 10.13662735161574
$ws.Cells.Item(24, 4).Value = 8.287845914984567
$ws.Cells.Item(24, 5).Value = 12.28289813263023
$ws.Cells.Item(24, 6).Value = 33.32913211468014
$ws.Cells.Item(24, 8).Value = 7.344005520526261
$ws.Cells.Item(24, 9).Value = 24.29144101399554
$ws.Cells.Item(24, 10).Value = 9.544466529482953
$ws.Cells.Item(24, 12).Value = 11.65553438624504
$ws.Cells.Item(24, 15).Value = 25.48312956556329
$ws.Cells.Item(25, 2).Value = 19.76972465227891
$ws.Cells.Item(25, 3).Value = 9.465311364035921
$ws.Cells.Item(25, 4).Value = 8.26774210946826
$ws.Cells.Item(25, 5).Value = 12.34052151164507
$ws.Cells.Item(25, 6).Value = 33.52264338063226
$ws.Cells.Item(25, 8).Value = 7.344005520526261
$ws.Cells.Item(25, 9).Value = 24.56194308863829
$ws.Cells.Item(25, 10).Value = 9.588580795133092
$ws.Cells.Item(25, 12).Value = 11.58396128793384
$ws.Cells.Item(25, 15).Value = 25.67317611036766
